$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H37").Value = -1

$ws.Range("G38").Value = 0.005
$ws.Range("H38").Value = 1

$ws.Range("H41").Value = -1

$ws.Range("H43").Value = -1

$ws.Range("H51").Value = -1

$ws.Range("H53").Value = -1
